# Insert a new price-record row at row 191 (this week's new Camote quote),
# pushing the existing historical rows 191-239 down to 192-240.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(191).Insert()

$ws.Cells.Item(191, 1).Value  = 10
$ws.Cells.Item(191, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(191, 3).Value  = "La Araucanía"
$ws.Cells.Item(191, 4).Value  = 45211
$ws.Cells.Item(191, 5).Value  = 9
$ws.Cells.Item(191, 6).Value  = 100114002
$ws.Cells.Item(191, 7).Value  = "Camote"
$ws.Cells.Item(191, 8).Value  = "Sin especificar"
$ws.Cells.Item(191, 9).Value  = "Primera"
$ws.Cells.Item(191, 10).Value = 140
$ws.Cells.Item(191, 11).Value = 24000
$ws.Cells.Item(191, 12).Value = 24000
$ws.Cells.Item(191, 13).Value = 24000
$ws.Cells.Item(191, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(191, 15).Value = "Perú"
$ws.Cells.Item(191, 16).Value = 1333
$ws.Cells.Item(191, 17).Value = 18
$ws.Cells.Item(191, 18).Value = "Hortaliza"
